$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Functional analyst with over 5 years of experience in support and enhancement of **SAP** systems (**MM**, **FI**, **PM**), focusing on implementation of improvements and functionalities, process automation, and incident management. Experienced in leading integration testing, working in **Agile** environments, and collaborating with global technical and functional teams. Capable of translating business requirements into clear technical specifications, managing end-to-end implementations."
$ws.Range("B6").Value = "Energy`nConsumer Goods & Services"
$ws.Range("B7").Value = "SAP FI configuration`nProcess automation`nIncident management`nIntegration testing`nFunctional design`nUser acceptance testing (UAT)`nData management"
$ws.Range("B8").Value = "First Certificate in English – FCE`nCaptton course SAP MM-SAD"
$ws.Range("B10").Value = "SAP FI Functional Analyst  `n**Gestión** y **configuración funcional** en **SAP FI**, incluyendo activos fijos y cuentas contables.  `n**Seguimiento** de incidencias y **análisis de errores** en procesos automatizados.  `n**Colaboración** con equipos técnicos y globales para la mejora continua del sistema."
$ws.Range("B11").Value = "Acquisition and Corporate SAP Analyst  `n**Liderazgo** en la implementación de mejoras y nuevas funcionalidades en **SAP**, incluyendo diseño funcional y coordinación de desarrollo.  `n**Conducción** de reuniones con áreas de negocio para el relevamiento de requerimientos y **optimización** de procesos.  `n**Diseño** y **ejecución** de mejoras en **SAP** y sus interfaces como **FIORI** y **ARIBA**."
$ws.Range("B12").Value = "Application Development Associate Junior  `n**Soporte** en **SAP MM**: resolución de incidentes y atención al usuario.  `n**Participación** en el proyecto de implementación de **SAP MM** y **PM**, incluyendo migración de datos en la industria petrolera."
$ws.Range("B13").Value = "IT Trainee  `n**Soporte** a la interfaz **SAP/GEP**: monitoreo de errores y **testing** de mejoras.  `n**Colaboración** en la mejora de procesos y soporte administrativo."
